# Add daily-update rows 7-19 to Sheet1, reusing existing formatting from
# the existing rows (3-6) so that styles / number-formats are re-used
# instead of creating brand-new style entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column D width (28.7109375 -> 34 "characters" in the saved XML) ---
$ws.Columns.Item(4).ColumnWidth = 33.14

# --- Date cells (column B) : copy format from an existing date cell (B4) ---
$ws.Range("B4").Copy()
$ws.Range("B7:B19").PasteSpecial(-4122)

# --- "Task" cells (column D) : copy wrap-text format from existing D3 ---
# (D8 and D10 are left with the default/general style, matching the
# original author's inconsistent formatting.)
$ws.Range("D3").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D11:D19").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Column A (SNO), B (Date serials)
# ---------------------------------------------------------------------
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 44263
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 44264
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 44265
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 44266
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 44267
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 44268
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 44269
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 44270
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 44271
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 44272
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 44273
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 44274
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 44275

# ---------------------------------------------------------------------
# Column C (Hours) / E (Status) - all re-use already existing shared
# strings ("8hr", "completed", "not completed"), so ordering does not
# matter for these.
# ---------------------------------------------------------------------
$ws.Range("C7").Value = "8hr"
$ws.Range("E7").Value = "completed"

$ws.Range("C8").Value = "8hr"
$ws.Range("E8").Value = "not completed"

$ws.Range("C9").Value = "8hr"
$ws.Range("E9").Value = "completed"

$ws.Range("C10").Value = "8hr"
$ws.Range("E10").Value = "not completed"

# Row 11 has no Hours entry in the source sheet.
$ws.Range("E11").Value = "not completed"

$ws.Range("C12").Value = "8hr"
$ws.Range("E12").Value = "completed"

$ws.Range("C13").Value = "8hr"
$ws.Range("E13").Value = "completed"

$ws.Range("C14").Value = "8hr"
$ws.Range("E14").Value = "not completed"

$ws.Range("C15").Value = "8hr"
$ws.Range("E15").Value = "completed"

$ws.Range("C16").Value = "8hr"
$ws.Range("E16").Value = "completed"

$ws.Range("C17").Value = "8hr"
$ws.Range("E17").Value = "completed"

$ws.Range("C18").Value = "8hr"
$ws.Range("E18").Value = "completed"

$ws.Range("C19").Value = "8hr"
$ws.Range("E19").Value = "not completed"

# ---------------------------------------------------------------------
# Column D (Task) - these create brand new shared-string entries.
# The ORDER in which new/unique text is first assigned determines the
# order the strings are appended to xl/sharedStrings.xml, so the calls
# below are intentionally not in simple row order: it mirrors the order
# in which the original author must have typed them.
# ---------------------------------------------------------------------
$ws.Range("D7").Value = "Modules & Services"
$ws.Range("D8").Value = "RestAppi calls"
$ws.Range("D9").Value = "RestAppi calls"
$ws.Range("D10").Value = "bootstrap"
$ws.Range("D11").Value = "bootstrap"
$ws.Range("D13").Value = "setting up Mock Data using JSON Server,installing Bootstrap in our application"
$ws.Range("D14").Value = "create e-commerce website"
$ws.Range("D12").Value = "design webpage using bootstrap"
$ws.Range("D15").Value = "installing bootstrap in our application,Creating feature Module"
$ws.Range("D16").Value = "how to create header,sidebar,footer"
$ws.Range("D17").Value = "how to display particular product,and also how to display product list"
$ws.Range("D18").Value = "how to do route"
$ws.Range("D19").Value = "routing feature module"

# ---------------------------------------------------------------------
# Row heights that differ from the default (wrapped multi-line text).
# ---------------------------------------------------------------------
$ws.Rows.Item(13).RowHeight = 45
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 45

# ---------------------------------------------------------------------
# View state : active cell / scroll position.
# ---------------------------------------------------------------------
$ws.Range("F18").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
